$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) added to the right of the existing "sum" column (G).
# Copy G1's formatting (bold header font, border, centered/top alignment)
# onto the new header cell so it reuses the same style as the rest of the
# header row, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New column's data values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
